$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:I6").Select()
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "3 trials both models"

$ws2.Range("A1").Value = "Type - Logical Augmentation LXMERT"
$ws2.Range("B1").Value = "Test acc."
$ws2.Range("C1").Value = "Dev acc."
$ws2.Range("D1").Value = "Min test"
$ws2.Range("E1").Value = "Max test"
$ws2.Range("F1").Value = "Mean test"
$ws2.Range("G1").Value = "Min dev"
$ws2.Range("H1").Value = "Max dev"
$ws2.Range("I1").Value = "Mean dev"
$ws2.Range("A2").Value = "Run#1"
$ws2.Range("B2").Value = 77.31
$ws2.Range("C2").Value = 75.77
$ws2.Range("A3").Value = "Run#2"
$ws2.Range("B3").Value = 77.91
$ws2.Range("C3").Value = 77.05
$ws2.Range("A4").Value = "Run#3"
$ws2.Range("B4").Value = 78.31
$ws2.Range("C4").Value = 77.05
$ws2.Range("D2").Formula = "=MIN(B2:B6)"
$ws2.Range("E2").Formula = "=MAX(B2:B6)"
$ws2.Range("F2").Formula = "=AVERAGE(B2:B6)"
$ws2.Range("G2").Formula = "=MIN(C2:C6)"
$ws2.Range("H2").Formula = "=MAX(C2:C6)"
$ws2.Range("I2").Formula = "=AVERAGE(C2:C6)"
$ws2.Range("F3").Formula = "=MAX(F2-D2,E2-F2)"
$ws2.Range("I3").Formula = "=MAX(H2-I2,I2-G2)"

$ws2.Range("A6").Value = "Type - Logical Augmentation VILT"
$ws2.Range("B6").Value = "Test acc."
$ws2.Range("C6").Value = "Dev acc."
$ws2.Range("D6").Value = "Min test"
$ws2.Range("E6").Value = "Max test"
$ws2.Range("F6").Value = "Mean test"
$ws2.Range("G6").Value = "Min dev"
$ws2.Range("H6").Value = "Max dev"
$ws2.Range("I6").Value = "Mean dev"
$ws2.Range("A7").Value = "Run#1"
$ws2.Range("B7").Value = 74.84
$ws2.Range("C7").Value = 73.29
$ws2.Range("A8").Value = "Run#2"
$ws2.Range("B8").Value = 73.84
$ws2.Range("C8").Value = 73.39
$ws2.Range("A9").Value = "Run#3"
$ws2.Range("B9").Value = 74.1
$ws2.Range("C9").Value = 73.39
$ws2.Range("D7").Formula = "=MIN(B7:B11)"
$ws2.Range("E7").Formula = "=MAX(B7:B11)"
$ws2.Range("F7").Formula = "=AVERAGE(B7:B11)"
$ws2.Range("G7").Formula = "=MIN(C7:C11)"
$ws2.Range("H7").Formula = "=MAX(C7:C11)"
$ws2.Range("I7").Formula = "=AVERAGE(C7:C11)"
$ws2.Range("F8").Formula = "=MAX(F7-D7,E7-F7)"
$ws2.Range("I8").Formula = "=MAX(H7-I7,I7-G7)"

$ws2.Range("A11").Value = "Type - Image Augmentation LXMERT"
$ws2.Range("B11").Value = "Test acc."
$ws2.Range("C11").Value = "Dev acc."
$ws2.Range("D11").Value = "Min test"
$ws2.Range("E11").Value = "Max test"
$ws2.Range("F11").Value = "Mean test"
$ws2.Range("G11").Value = "Min dev"
$ws2.Range("H11").Value = "Max dev"
$ws2.Range("I11").Value = "Mean dev"
$ws2.Range("A12").Value = "Run#1"
$ws2.Range("B12").Value = 75.28
$ws2.Range("C12").Value = 71.81
$ws2.Range("A13").Value = "Run#2"
$ws2.Range("B13").Value = 76.42
$ws2.Range("C13").Value = 72.9
$ws2.Range("A14").Value = "Run#3"
$ws2.Range("B14").Value = 76.57
$ws2.Range("C14").Value = 73.89
$ws2.Range("D12").Formula = "=MIN(B12:B16)"
$ws2.Range("E12").Formula = "=MAX(B12:B16)"
$ws2.Range("F12").Formula = "=AVERAGE(B12:B16)"
$ws2.Range("G12").Formula = "=MIN(C12:C16)"
$ws2.Range("H12").Formula = "=MAX(C12:C16)"
$ws2.Range("I12").Formula = "=AVERAGE(C12:C16)"
$ws2.Range("F13").Formula = "=MAX(F12-D12,E12-F12)"
$ws2.Range("I13").Formula = "=MAX(H12-I12,I12-G12)"

$ws2.Range("A16").Value = "Type - Image Augmentation VILT"
$ws2.Range("B16").Value = "Test acc."
$ws2.Range("C16").Value = "Dev acc."
$ws2.Range("D16").Value = "Min test"
$ws2.Range("E16").Value = "Max test"
$ws2.Range("F16").Value = "Mean test"
$ws2.Range("G16").Value = "Min dev"
$ws2.Range("H16").Value = "Max dev"
$ws2.Range("I16").Value = "Mean dev"
$ws2.Range("A17").Value = "Run#1"
$ws2.Range("A18").Value = "Run#2"
$ws2.Range("A19").Value = "Run#3"
$ws2.Range("D17").Formula = "=MIN(B17:B21)"
$ws2.Range("E17").Formula = "=MAX(B17:B21)"
$ws2.Range("F17").Formula = "=AVERAGE(B17:B21)"
$ws2.Range("G17").Formula = "=MIN(C17:C21)"
$ws2.Range("H17").Formula = "=MAX(C17:C21)"
$ws2.Range("I17").Formula = "=AVERAGE(C17:C21)"
$ws2.Range("F18").Formula = "=MAX(F17-D17,E17-F17)"
$ws2.Range("I18").Formula = "=MAX(H17-I17,I17-G17)"

$ws2.Range("A21").Value = "Type - Contrastive LXMERT"
$ws2.Range("B21").Value = "Test acc."
$ws2.Range("C21").Value = "Dev acc."
$ws2.Range("D21").Value = "Min test"
$ws2.Range("E21").Value = "Max test"
$ws2.Range("F21").Value = "Mean test"
$ws2.Range("G21").Value = "Min dev"
$ws2.Range("H21").Value = "Max dev"
$ws2.Range("I21").Value = "Mean dev"
$ws2.Range("A22").Value = "Run#1"
$ws2.Range("B22").Value = 77.46
$ws2.Range("C22").Value = 75.67
$ws2.Range("A23").Value = "Run#2"
$ws2.Range("B23").Value = 79
$ws2.Range("C23").Value = 76.85
$ws2.Range("A24").Value = "Run#3"
$ws2.Range("B24").Value = 77.27
$ws2.Range("C24").Value = 76.46
$ws2.Range("D22").Formula = "=MIN(B22:B26)"
$ws2.Range("E22").Formula = "=MAX(B22:B26)"
$ws2.Range("F22").Formula = "=AVERAGE(B22:B26)"
$ws2.Range("G22").Formula = "=MIN(C22:C26)"
$ws2.Range("H22").Formula = "=MAX(C22:C26)"
$ws2.Range("I22").Formula = "=AVERAGE(C22:C26)"
$ws2.Range("F23").Formula = "=MAX(F22-D22,E22-F22)"
$ws2.Range("I23").Formula = "=MAX(H22-I22,I22-G22)"

$ws2.Range("A26").Value = "Type - Contrastive VILT"
$ws2.Range("B26").Value = "Test acc."
$ws2.Range("C26").Value = "Dev acc."
$ws2.Range("D26").Value = "Min test"
$ws2.Range("E26").Value = "Max test"
$ws2.Range("F26").Value = "Mean test"
$ws2.Range("G26").Value = "Min dev"
$ws2.Range("H26").Value = "Max dev"
$ws2.Range("I26").Value = "Mean dev"
$ws2.Range("A27").Value = "Run#1"
$ws2.Range("A28").Value = "Run#2"
$ws2.Range("A29").Value = "Run#3"
$ws2.Range("D27").Formula = "=MIN(B27:B31)"
$ws2.Range("E27").Formula = "=MAX(B27:B31)"
$ws2.Range("F27").Formula = "=AVERAGE(B27:B31)"
$ws2.Range("G27").Formula = "=MIN(C27:C31)"
$ws2.Range("H27").Formula = "=MAX(C27:C31)"
$ws2.Range("I27").Formula = "=AVERAGE(C27:C31)"
$ws2.Range("F28").Formula = "=MAX(F27-D27,E27-F27)"
$ws2.Range("I28").Formula = "=MAX(H27-I27,I27-G27)"

$ws2.Range("K21").Value = "Type - Contrastive LXMERT Trial 2"
$ws2.Range("L21").Value = "Test acc."
$ws2.Range("M21").Value = "Dev acc."
$ws2.Range("N21").Value = "Min test"
$ws2.Range("O21").Value = "Max test"
$ws2.Range("P21").Value = "Mean test"
$ws2.Range("Q21").Value = "Min dev"
$ws2.Range("R21").Value = "Max dev"
$ws2.Range("S21").Value = "Mean dev"
$ws2.Range("K22").Value = "Run#1"
$ws2.Range("L22").Value = 0.785042100049529
$ws2.Range("M22").Value = 0.7675568743818
$ws2.Range("K23").Value = "Run#2"
$ws2.Range("L23").Value = 0.779593858345715
$ws2.Range("M23").Value = 0.762611275964391
$ws2.Range("K24").Value = "Run#3"
$ws2.Range("L24").Value = 0.775136206042595
$ws2.Range("M24").Value = 0.751730959446093
$ws2.Range("N22").Formula = "=MIN(L22:L26)"
$ws2.Range("O22").Formula = "=MAX(L22:L26)"
$ws2.Range("P22").Formula = "=AVERAGE(L22:L26)"
$ws2.Range("Q22").Formula = "=MIN(M22:M26)"
$ws2.Range("R22").Formula = "=MAX(M22:M26)"
$ws2.Range("S22").Formula = "=AVERAGE(M22:M26)"
$ws2.Range("P23").Formula = "=MAX(P22-N22,O22-P22)"
$ws2.Range("S23").Formula = "=MAX(R22-S22,S22-Q22)"

$ws2.Range("A31").Value = "Type - SNLI-VE Pretraining LXMERT"
$ws2.Range("B31").Value = "Test acc."
$ws2.Range("C31").Value = "Dev acc."
$ws2.Range("D31").Value = "Min test"
$ws2.Range("E31").Value = "Max test"
$ws2.Range("F31").Value = "Mean test"
$ws2.Range("G31").Value = "Min dev"
$ws2.Range("H31").Value = "Max dev"
$ws2.Range("I31").Value = "Mean dev"
$ws2.Range("A32").Value = "Run#1"
$ws2.Range("B32").Value = 75.45
$ws2.Range("C32").Value = 75.88
$ws2.Range("A33").Value = "Run#2"
$ws2.Range("B33").Value = 75.6
$ws2.Range("C33").Value = 75.98
$ws2.Range("A34").Value = "Run#3"
$ws2.Range("B34").Value = 75.29
$ws2.Range("C34").Value = 76.12
$ws2.Range("D32").Formula = "=MIN(B32:B36)"
$ws2.Range("E32").Formula = "=MAX(B32:B36)"
$ws2.Range("F32").Formula = "=AVERAGE(B32:B36)"
$ws2.Range("G32").Formula = "=MIN(C32:C36)"
$ws2.Range("H32").Formula = "=MAX(C32:C36)"
$ws2.Range("I32").Formula = "=AVERAGE(C32:C36)"
$ws2.Range("F33").Formula = "=MAX(F32-D32,E32-F32)"
$ws2.Range("I33").Formula = "=MAX(H32-I32,I32-G32)"

$ws2.Range("A36").Value = "Type - SNLI-VE Finetuning LXMERT"
$ws2.Range("B36").Value = "Test acc."
$ws2.Range("C36").Value = "Dev acc."
$ws2.Range("D36").Value = "Min test"
$ws2.Range("E36").Value = "Max test"
$ws2.Range("F36").Value = "Mean test"
$ws2.Range("G36").Value = "Min dev"
$ws2.Range("H36").Value = "Max dev"
$ws2.Range("I36").Value = "Mean dev"
$ws2.Range("A37").Value = "Run#1"
$ws2.Range("B37").Value = 75.04
$ws2.Range("C37").Value = 76.24
$ws2.Range("A38").Value = "Run#2"
$ws2.Range("B38").Value = 74.84
$ws2.Range("C38").Value = 75.9
$ws2.Range("A39").Value = "Run#3"
$ws2.Range("B39").Value = 75.12
$ws2.Range("C39").Value = 76.2
$ws2.Range("D37").Formula = "=MIN(B37:B41)"
$ws2.Range("E37").Formula = "=MAX(B37:B41)"
$ws2.Range("F37").Formula = "=AVERAGE(B37:B41)"
$ws2.Range("G37").Formula = "=MIN(C37:C41)"
$ws2.Range("H37").Formula = "=MAX(C37:C41)"
$ws2.Range("I37").Formula = "=AVERAGE(C37:C41)"
$ws2.Range("F38").Formula = "=MAX(F37-D37,E37-F37)"
$ws2.Range("I38").Formula = "=MAX(H37-I37,I37-G37)"

$ws2.Columns("A").ColumnWidth = 53
$ws2.Columns("K").ColumnWidth = 33.5

$ws2.Activate()
$ws2.Range("E33").Select()